$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: add new date column header C1 (13-01-2023), matching B1's style ---
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value = "13-01-2023"

# --- Make sure column A style (bold / bordered) extends to the new row 9 ---
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122)

# --- Rewrite data rows 2-9 in their new order / with new values ---
# Row 2: Adcap Wise Capital Growth
$ws.Range("A2").Value = "Adcap Wise Capital Growth"
$ws.Range("B2").Value = 20133706.38
$ws.Range("C2").Value = 17244045.31

# Row 3: Alpha planeam equil
$ws.Range("A3").Value = "Alpha planeam equil"
$ws.Range("B3").Value = 157392.29
$ws.Range("C3").Value = 134611.49

# Row 4: Alpha renta balan global (new fund, no value for the first period)
$ws.Range("A4").Value = "Alpha renta balan global"
$ws.Range("B4").Value = ""
$ws.Range("C4").Value = 990698.95

# Row 5: Arpenta acciones
$ws.Range("A5").Value = "Arpenta acciones"
$ws.Range("B5").Value = 6946.92
$ws.Range("C5").Value = 6938.89

# Row 6: Delta Recursos Naturales
$ws.Range("A6").Value = "Delta Recursos Naturales"
$ws.Range("B6").Value = 347109.59
$ws.Range("C6").Value = 347581.5

# Row 7: HF Acciones Argentinas
$ws.Range("A7").Value = "HF Acciones Argentinas"
$ws.Range("B7").Value = 35692.54
$ws.Range("C7").Value = 31152.26

# Row 8: avg
$ws.Range("A8").Value = "avg"
$ws.Range("B8").Value = 4136169.54
$ws.Range("C8").Value = 3125838.07

# Row 9: total
$ws.Range("A9").Value = "total"
$ws.Range("B9").Value = 20680847.72
$ws.Range("C9").Value = 18755028.4
